{"js": "// Office.js (Word JavaScript API) edit script.\n// Splits two long, run-embedded sentences into multiple line breaks\n// (<w:br/>) within the same run, matching the target diff:\n//   1. The \"Crit\u00e9rio:\" evaluation rule paragraph gets 3 new line breaks.\n//   2. The \"Bibliografia\" paragraph gets 2 new line breaks.\n//\n// NOTE: this runtime's `insertBreak()` does not honor InsertLocation on a\n// search-result range (it always appends at the end of the containing\n// paragraph), so instead we insert the literal vertical-tab character\n// (U+000B) via `insertText(..., InsertLocation.before)`, which Word's OOXML\n// writer serializes as a proper <w:br/> element split across <w:t> runs -\n// exactly the shape the diff expects.\n\nasync function insertLineBreakBefore(scopeRange, searchText) {\n  const results = scopeRange.search(searchText, { matchCase: true, ignorePunct: false, ignoreSpace: false });\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + searchText);\n  }\n  results.items[0].insertText(\"\\v\", Word.InsertLocation.before);\n  await context.sync();\n}\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// --- 1. \"Crit\u00e9rio:\" paragraph (Avalia\u00e7\u00e3o section) ---------------------\nlet criterioParagraph = null;\nfor (const p of paragraphs.items) {\n  if (p.text.indexOf(\"MS= (2xP1+3xP2)/5\") !== -1) {\n    criterioParagraph = p;\n    break;\n  }\n}\nif (!criterioParagraph) {\n  throw new Error(\"Could not locate the Crit\u00e9rio paragraph.\");\n}\n\nawait insertLineBreakBefore(criterioParagraph, \"MS> ou = 5,0: Aluno Aprovado\");\nawait insertLineBreakBefore(criterioParagraph, \"MS< 3,0: Aluno Reprovado\");\nawait insertLineBreakBefore(criterioParagraph, \"3,0 < ou = MS < 5,0: Aluno de Recupera\u00e7\u00e3o.\");\n\n// --- 2. \"Bibliografia\" paragraph --------------------------------------\nlet bibliografiaParagraph = null;\nfor (const p of paragraphs.items) {\n  if (p.text.indexOf(\"McMURRY, J. Qu\u00edmica Org\u00e2nica\") !== -1) {\n    bibliografiaParagraph = p;\n    break;\n  }\n}\nif (!bibliografiaParagraph) {\n  throw new Error(\"Could not locate the Bibliografia paragraph.\");\n}\n\nawait insertLineBreakBefore(bibliografiaParagraph, \"- MORRISON, R.T. e BOYD, R.N.\");\nawait insertLineBreakBefore(bibliografiaParagraph, \"- SOLOMONS, T.W.G., FRYHLE, C.B.\");\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Splits two long, run-embedded sentences into multiple manual line breaks\n# (<w:br/>), matching the target diff:\n#   1. The \"Crit\u00e9rio:\" evaluation-rule paragraph gets 3 new line breaks.\n#   2. The \"Bibliografia\" paragraph gets 2 new line breaks.\n#\n# Approach: use Find/Replace on $d.Content with MatchWildcards = $false (so\n# the punctuation-heavy search strings are treated literally) and a\n# replacement string prefixed with \"^l\" - Word's special code for a manual\n# line break - which round-trips to an actual <w:br/> OOXML element.\n\n$d = $word.ActiveDocument\n\nfunction Insert-LineBreakBefore([string]$searchText) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $replaceText = \"^l\" + $searchText\n    $ok = $rng.Find.Execute(\n        $searchText,   # FindText\n        $true,         # MatchCase\n        $false,        # MatchWholeWord\n        $false,        # MatchWildcards\n        $false,        # MatchSoundsLike\n        $false,        # MatchAllWordForms\n        $true,         # Forward\n        1,             # Wrap (wdFindContinue)\n        $false,        # Format\n        $replaceText,  # ReplaceWith\n        2              # Replace (wdReplaceOne)\n    )\n    if (-not $ok) {\n        throw \"Text not found: $searchText\"\n    }\n}\n\n# --- 1. \"Crit\u00e9rio:\" paragraph (Avalia\u00e7\u00e3o section) ----------------------\nInsert-LineBreakBefore(\"MS> ou = 5,0: Aluno Aprovado\")\nInsert-LineBreakBefore(\"MS< 3,0: Aluno Reprovado\")\nInsert-LineBreakBefore(\"3,0 < ou = MS < 5,0: Aluno de Recupera\u00e7\u00e3o.\")\n\n# --- 2. \"Bibliografia\" paragraph ---------------------------------------\nInsert-LineBreakBefore(\"- MORRISON, R.T. e BOYD, R.N.\")\nInsert-LineBreakBefore(\"- SOLOMONS, T.W.G., FRYHLE, C.B.\")\n"}
